$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to text
# so Excel does not reinterpret them (matching the source data which
# stores every value as a string).
$textCells = @('D5','D6','D7','D9','D11','D12','D13','D14','D16','D19','D22','D23','D24','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D38','D39','D43','D44','D46','D47','D49','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the refreshed cryptos feed.
$ws.Range('D2').Value = '70.768.67'
$ws.Range('E2').Value = '  +5.28%  '
$ws.Range('D3').Value = '3.654.59'
$ws.Range('E3').Value = '  +5.25%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '592.88'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').Value = '194.77'
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('D7').Value = '0.649'
$ws.Range('E7').Value = '  +2.50%  '
$ws.Range('D8').Value = '3.647.51'
$ws.Range('E8').Value = '  +5.20%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +4.99%  '
$ws.Range('D11').Value = '0.676'
$ws.Range('E11').Value = '  +4.07%  '
$ws.Range('D12').Value = '58.62'
$ws.Range('E12').Value = '  +2.95%  '
$ws.Range('D13').Value = '0.0000294'
$ws.Range('E13').Value = '  +5.05%  '
$ws.Range('D14').Value = '9.98'
$ws.Range('E14').Value = '  +5.54%  '
$ws.Range('D15').Value = '4.235.34'
$ws.Range('E15').Value = '  +5.18%  '
$ws.Range('D16').Value = '20.00'
$ws.Range('E16').Value = '  +6.36%  '
$ws.Range('D17').Value = '3.649.27'
$ws.Range('E17').Value = '  +5.20%  '
$ws.Range('D18').Value = '70.777.00'
$ws.Range('E18').Value = '  +5.26%  '
$ws.Range('D19').Value = '12.84'
$ws.Range('E19').Value = '  +5.35%  '
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('E21').Value = '  +5.25%  '
$ws.Range('D22').Value = '492.40'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').Value = '18.85'
$ws.Range('E23').Value = '  +11.86%  '
$ws.Range('D24').Value = '5.37'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').Value = '4.51'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').Value = '91.84'
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('D27').Value = '3.17'
$ws.Range('E27').Value = '  +7.45%  '
$ws.Range('D28').Value = '11.51'
$ws.Range('E28').Value = '  +4.86%  '
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '33.11'
$ws.Range('E30').Value = '  +5.30%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '7.92'
$ws.Range('E31').Value = '  +10.07%  '
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +8.59%  '
$ws.Range('D33').Value = '631.47'
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('D34').Value = '12.33'
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('D35').Value = '65.65'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('D36').Value = '40.66'
$ws.Range('E36').Value = '  +10.87%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0840'
$ws.Range('E37').Value = '  +10.49%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').Value = '0.414'
$ws.Range('E38').Value = '  +7.39%  '
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('D42').Value = '3.317.59'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').Value = '3.18'
$ws.Range('E43').Value = '  +9.22%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.87'
$ws.Range('E44').Value = '  +14.13%  '
$ws.Range('E45').Value = '  +5.83%  '
$ws.Range('D46').Value = '2.92'
$ws.Range('E46').Value = '  +5.27%  '
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('D49').Value = '9.26'
$ws.Range('E49').Value = '  +5.35%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  -0.19%  '

# Restore default (General) styling on the cells we forced to text so the
# workbook formatting matches the original (text-typed) cells exactly.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
